$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update confirmed-in-hospital (治愈/cured column C) values:
# 罗田 (row 6): 2 -> 0
# 麻城 (row 11): 0 -> 2
$ws.Range("C6").Value = 0
$ws.Range("C11").Value = 2

# Update selection to reflect new active cell
$ws.Range("F11").Select()
